# Sprint 43 test case report: fill in the "Day 3" test case summary block
# (rows 15-17) with the day's written / execution / review counts, and
# leave the selection on the last cell that was entered (C17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = 7040
$ws.Range("C16").Value = 2338
$ws.Range("C17").Value = 2338

$ws.Range("C17").Select()
